# Shift the "date" column forward by 3 days (2025-04-28..2025-05-05 ->
# 2025-05-01..2025-05-08) and update the new "history" counts that came
# in with the latest days (B6=26, B8=2, B9=4). Column A holds its dates
# as plain text (t="str") in the source file, so we must avoid letting
# Excel auto-convert the "YYYY-MM-DD" strings into real date serials
# (which is what a plain `.Value = "2025-05-01"` assignment would do).
#
# Trick: write the text via a formula that evaluates to a string
# (="2025-05-01"), which Excel always treats as text regardless of its
# look, then copy/paste-special-values it back onto itself. That bakes
# in a literal text value with no residual formula and no style churn.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @(
    "2025-05-01",
    "2025-05-02",
    "2025-05-03",
    "2025-05-04",
    "2025-05-05",
    "2025-05-06",
    "2025-05-07",
    "2025-05-08"
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 1)
    $cell.Formula = "=""" + $dates[$i] + """"
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
}

$excel.CutCopyMode = 0

# New history values for the added days.
$ws.Range("B6").Value = 26
$ws.Range("B8").Value = 2
$ws.Range("B9").Value = 4
